$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Remove the 7 detail rows for "HỆ THỐNG" (rows 4-10)
$ws.Rows("4:10").Delete()

# After the shift, "Tổng lương tại HỆ THỐNG" (originally row 35) is now at row 28
$ws.Rows("28:28").Delete()
